$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "Averaged Data" placeholder values ---
# Row 2 (Proof of Work) gets a real TPS number; every other new cell in the
# table body is marked "N/A" (data not available / not yet collected).

# B2 is numeric (TPS for Proof of Work)
$ws.Range("B2").Value = 17

# Remaining cells in row 2
$ws.Range("C2:G2").Value = "N/A"

# Rows 3-10, columns B:G all become "N/A"
$ws.Range("B3:G10").Value = "N/A"

# --- Number formats per column, matching the table's column intent ---
# % of nodes required to take over network -> percent (already was)
$ws.Range("E2:E10").NumberFormat = "0.00%"
# Strengths / Weaknesses -> plain text format
$ws.Range("F2:G10").NumberFormat = "@"
# TPS / Energy Use / Nakamoto Coefficient -> 2-decimal numeric format
$ws.Range("B2:D10").NumberFormat = "0.00"

# --- Misc view state ---
$ws.Range("E18").Select()
